$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "53.874.76"
$ws.Range("E2").Value = "  -4.28%  "

# Row 3
$ws.Range("D3").Value = "2.237.51"
$ws.Range("E3").Value = "  -5.76%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'488.17"
$ws.Range("E5").Value = "  -2.70%  "

# Row 6
$ws.Range("D6").Value = "'126.64"
$ws.Range("E6").Value = "  -2.77%  "

# Row 7
$ws.Range("E7").Value = "  +0.28%  "

# Row 8
$ws.Range("E8").Value = "  -4.56%  "

# Row 9
$ws.Range("D9").Value = "2.247.19"
$ws.Range("E9").Value = "  -5.53%  "

# Row 10
$ws.Range("D10").Value = "'0.0916"
$ws.Range("E10").Value = "  -6.84%  "

# Row 12
$ws.Range("E12").Value = "  +0.85%  "

# Row 13
$ws.Range("E13").Value = "  -2.87%  "

# Row 14
$ws.Range("D14").Value = "2.635.69"
$ws.Range("E14").Value = "  -5.70%  "

# Row 15
$ws.Range("D15").Value = "'21.15"
$ws.Range("E15").Value = "  -2.11%  "

# Row 16
$ws.Range("D16").Value = "53.788.71"
$ws.Range("E16").Value = "  -4.33%  "

# Row 17
$ws.Range("E17").Value = "  -3.11%  "

# Row 18
$ws.Range("D18").Value = "2.252.19"
$ws.Range("E18").Value = "  -6.45%  "

# Row 19
$ws.Range("E19").Value = "  -1.05%  "

# Row 20
$ws.Range("D20").Value = "'9.59"
$ws.Range("E20").Value = "  -4.68%  "

# Row 21
$ws.Range("D21").Value = "'299.81"
$ws.Range("E21").Value = "  -2.37%  "

# Row 22
$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  -2.49%  "

# Row 23
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.41%  "

# Row 24
$ws.Range("D24").Value = "'63.66"
$ws.Range("E24").Value = "  -1.44%  "

# Row 26
$ws.Range("E26").Value = "  -1.01%  "

# Row 27
$ws.Range("E27").Value = "  -2.75%  "

# Row 28
$ws.Range("D28").Value = "'7.07"
$ws.Range("E28").Value = "  -3.43%  "

# Row 29
$ws.Range("D29").Value = "'169.15"
$ws.Range("E29").Value = "  -1.82%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0688"
$ws.Range("E30").Value = "  -3.66%  "

# Row 31
$ws.Range("E31").Value = "  -2.58%  "

# Row 33
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'5.73"
$ws.Range("E33").Value = "  -0.56%  "

# Row 34
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "  -0.05%  "

# Row 35
$ws.Range("E35").Value = "  -2.27%  "

# Row 36
$ws.Range("D36").Value = "'17.42"
$ws.Range("E36").Value = "  -0.81%  "

# Row 37
$ws.Range("E37").Value = "  -0.94%  "

# Row 38
$ws.Range("E38").Value = "  +6.88%  "

# Row 39
$ws.Range("E39").Value = "  -5.09%  "

# Row 40
$ws.Range("D40").Value = "'35.67"
$ws.Range("E40").Value = "  -0.91%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.37"
$ws.Range("E41").Value = "  -2.27%  "

# Row 42
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.366"
$ws.Range("E42").Value = "  -0.80%  "

# Row 43
$ws.Range("E43").Value = "  -1.45%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'4.81"
$ws.Range("E44").Value = "  +1.28%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'122.34"
$ws.Range("E45").Value = "  -6.48%  "

# Row 46
$ws.Range("E46").Value = "  -2.87%  "

# Row 47
$ws.Range("E47").Value = "  -5.59%  "

# Row 48
$ws.Range("D48").Value = "'236.66"
$ws.Range("E48").Value = "  -1.95%  "

# Row 49
$ws.Range("D49").Value = "'0.0471"
$ws.Range("E49").Value = "  -2.63%  "

# Row 50
$ws.Range("D50").Value = "'0.0203"
$ws.Range("E50").Value = "  -2.93%  "

# Row 51
$ws.Range("D51").Value = "'16.16"
$ws.Range("E51").Value = "  -4.17%  "
